$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the old row 33 ("2_state" header group)
$ws.Rows("33:35").Insert()

# Change CREDENTIAL_TYPE value from "window" to "asset"
$ws.Range("C22").Value = "asset"

# New keys in column A
$ws.Range("A33").Value = "TIME_CHECK_IN"
$ws.Range("A34").Value = "TIME_CHECK_OUT"
$ws.Range("A35").Value = "LAST_TIME_CHECK_OUT"

# New descriptions in column E (written bottom-up)
$ws.Range("E35").Value = "giờ tan làm ngày thứ 7"
$ws.Range("E34").Value = "tan làm ngày bình thường"
$ws.Range("E33").Value = "giờ vào làm"

# New time values in column C
$ws.Range("C33").Value = 0.35416666666666669
$ws.Range("C34").Value = 0.72916666666666663
$ws.Range("C35").Value = 0.5

$ws.Range("C33:C35").NumberFormat = "h:mm"

# Apply consistent font across the new rows
$ws.Range("A33:E35").Font.Name = "Open Sans"
$ws.Range("A33:E35").Font.Size = 11
